$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '27.872.75'
Set-TextValue 'D3' '1.768.86'
$ws.Range('E3').Value = '  +0.94%  '
Set-TextValue 'D4' '1.001'
$ws.Range('E4').Value = '  -0.21%  '
Set-TextValue 'D5' '327.43'
$ws.Range('E5').Value = '  +1.04%  '
$ws.Range('E6').Value = '  -0.22%  '
Set-TextValue 'D7' '0.4481'
$ws.Range('E7').Value = '  -2.49%  '
Set-TextValue 'D8' '0.3551'
$ws.Range('E8').Value = '  -0.70%  '
Set-TextValue 'D9' '0.07457'
$ws.Range('E9').Value = '  -0.49%  '
Set-TextValue 'D10' '42.06'
$ws.Range('E10').Value = '  -0.18%  '
Set-TextValue 'D11' '1.096'
$ws.Range('E11').Value = '  +0.16%  '
Set-TextValue 'D12' '1.000'
$ws.Range('E12').Value = '  -0.22%  '
Set-TextValue 'D13' '20.86'
$ws.Range('E13').Value = '  +0.53%  '
Set-TextValue 'D14' '6.027'
$ws.Range('E14').Value = '  +0.33%  '
Set-TextValue 'D15' '7.207'
$ws.Range('E15').Value = '  +1.71%  '
Set-TextValue 'D16' '1.768.86'
$ws.Range('E16').Value = '  +0.92%  '
Set-TextValue 'D17' '93.22'
$ws.Range('E17').Value = '  +0.91%  '
Set-TextValue 'D18' '0.00001058'
$ws.Range('E18').Value = '  -0.65%  '
Set-TextValue 'D19' '0.06431'
$ws.Range('E19').Value = '  +0.02%  '
Set-TextValue 'D20' '0.9998'
$ws.Range('E20').Value = '  -0.27%  '
$ws.Range('E21').Value = '  +2.38%  '
$ws.Range('E22').Value = '  -0.63%  '
Set-TextValue 'D23' '27.912.39'
$ws.Range('E23').Value = '  +0.93%  '
Set-TextValue 'D24' '11.29'
$ws.Range('E24').Value = '  +0.39%  '
Set-TextValue 'D25' '2.109'
$ws.Range('E25').Value = '  +0.05%  '
Set-TextValue 'D26' '161.72'
$ws.Range('E26').Value = '  -1.65%  '
Set-TextValue 'D27' '20.24'
$ws.Range('E27').Value = '  -0.36%  '
Set-TextValue 'D28' '1.970.88'
$ws.Range('E28').Value = '  +0.83%  '
Set-TextValue 'D29' '2.166'
$ws.Range('E29').Value = '  +4.39%  '
Set-TextValue 'D30' '125.05'
$ws.Range('E30').Value = '  -1.20%  '
Set-TextValue 'D31' '1.096'
$ws.Range('E31').Value = '  +3.21%  '
Set-TextValue 'D32' '0.09161'
$ws.Range('E32').Value = '  -0.20%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D33' '5.581'
$ws.Range('E33').Value = '  +0.85%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D34' '3.646'
$ws.Range('E34').Value = '  -0.47%  '
Set-TextValue 'D35' '11.90'
$ws.Range('E35').Value = '  +0.07%  '
Set-TextValue 'D36' '0.02294'
$ws.Range('E36').Value = '  -0.01%  '
Set-TextValue 'D37' '0.06114'
$ws.Range('E37').Value = '  +1.00%  '
Set-TextValue 'D38' '0.2096'
$ws.Range('E38').Value = '  -0.20%  '
Set-TextValue 'D39' '4.963'
$ws.Range('E39').Value = '  -0.33%  '
Set-TextValue 'D40' '0.6299'
$ws.Range('E40').Value = '  -0.45%  '
Set-TextValue 'D41' '1.182'
$ws.Range('E41').Value = '  -2.21%  '
Set-TextValue 'D42' '1.393'
$ws.Range('E42').Value = '  +1.01%  '
Set-TextValue 'D43' '7.951'
$ws.Range('E43').Value = '  +2.22%  '
Set-TextValue 'D44' '13.25'
$ws.Range('E44').Value = '  +0.19%  '
Set-TextValue 'D45' '3.744'
$ws.Range('E45').Value = '  +0.87%  '
Set-TextValue 'D46' '0.5864'
$ws.Range('E46').Value = '  -0.86%  '
Set-TextValue 'D47' '122.43'
$ws.Range('E47').Value = '  -0.51%  '
Set-TextValue 'D48' '1.950'
$ws.Range('E48').Value = '  +0.41%  '
Set-TextValue 'D49' '0.06902'
$ws.Range('E49').Value = '  +0.62%  '
Set-TextValue 'D50' '1.135'
$ws.Range('E50').Value = '  -0.60%  '
Set-TextValue 'D51' '72.88'
$ws.Range('E51').Value = '  +1.03%  '
